$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of course data (dates/booleans are entered with a leading
# apostrophe so Excel stores them as literal text instead of auto-converting
# them to date serials / real booleans).
$rows = @(
    @{ Row = 3; Date = "03/24/2019"; Count = 11 },
    @{ Row = 4; Date = "03/24/2020"; Count = 12 },
    @{ Row = 5; Date = "03/24/2021"; Count = 14 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row`:D$row").HorizontalAlignment = -4131

    $ws.Range("A$row").Value = "'" + $r.Date
    $ws.Range("B$row").Value = $r.Count
    $ws.Range("C$row").Value = "'true"
    $ws.Range("D$row").Value = "Chemistry"
}

$ws.Range("B11").Select() | Out-Null
